# Auto-generated Excel COM-interop script to apply crypto price/volume updates
# Uses a formula+copy+paste-special(values) trick to force text-typed cell values
# (avoids Excel auto-converting numeric-looking strings like "237.98" into floats,
#  and avoids introducing any new cell styles / number formats).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $escaped = $val -replace '"', '""'
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

$updates = [ordered]@{
    "D2" = '29.114.92'
    "E2" = '  -2.35%  '
    "D3" = '1.851.46'
    "E3" = '  -1.13%  '
    "D4" = '0.9980'
    "E4" = '  -0.20%  '
    "D5" = '0.6955'
    "E5" = '  -4.93%  '
    "D6" = '237.98'
    "D7" = '0.9991'
    "E7" = '  -0.23%  '
    "D8" = '0.07740'
    "E8" = '  +8.56%  '
    "E9" = '  -3.00%  '
    "D10" = '23.33'
    "E10" = '  -4.41%  '
    "D11" = '0.08124'
    "E11" = '  -0.37%  '
    "D12" = '1.872.59'
    "E12" = '  -0.77%  '
    "D13" = '0.7266'
    "E13" = '  -2.06%  '
    "D14" = '5.213'
    "E14" = '  -2.41%  '
    "D15" = '89.04'
    "D16" = '29.108.99'
    "B17" = 'ShibaInu'
    "C17" = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    "D17" = '0.000007843'
    "E17" = '  +0.46%  '
    "B18" = 'Uniswap'
    "C18" = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    "D18" = '5.754'
    "E18" = '  -4.17%  '
    "D19" = '13.21'
    "E19" = '  -1.30%  '
    "D20" = '236.29'
    "E20" = '  -4.88%  '
    "E21" = '  -0.33%  '
    "D22" = '2.099.19'
    "E22" = '  -2.80%  '
    "D23" = '0.9985'
    "E23" = '  -0.16%  '
    "D24" = '7.608'
    "E24" = '  -1.87%  '
    "D25" = '8.984'
    "E25" = '  -2.45%  '
    "D26" = '161.16'
    "D27" = '0.1434'
    "E27" = '  -6.91%  '
    "D28" = '18.08'
    "E28" = '  -2.44%  '
    "E29" = '  -1.88%  '
    "D30" = '1.397'
    "E30" = '  -3.58%  '
    "D31" = '4.499'
    "E31" = '  -0.51%  '
    "D32" = '1.488'
    "E32" = '  -2.22%  '
    "D33" = '4.012'
    "E33" = '  -4.05%  '
    "D34" = '0.05238'
    "E34" = '  -1.39%  '
    "D35" = '1.183'
    "E35" = '  -3.88%  '
    "D36" = '0.7065'
    "E36" = '  -4.61%  '
    "D37" = '1.024'
    "E37" = '  +2.26%  '
    "D38" = '2.646'
    "E38" = '  -2.15%  '
    "D39" = '0.01856'
    "E39" = '  -4.21%  '
    "D40" = '2.671'
    "E40" = '  -2.33%  '
    "D41" = '0.9141'
    "E41" = '  +5.18%  '
    "D42" = '1.095.30'
    "E42" = '  +4.54%  '
    "D43" = '5.985'
    "E43" = '  +0.35%  '
    "D44" = '0.4275'
    "E44" = '  -4.24%  '
    "D45" = '70.82'
    "E45" = '  -0.59%  '
    "D46" = '0.9991'
    "E46" = '  -0.22%  '
    "D47" = '102.99'
    "E47" = '  -0.90%  '
    "D48" = '1.771'
    "E48" = '  -2.53%  '
    "D49" = '1.994.93'
    "E49" = '  -2.90%  '
    "D50" = '9.184'
    "E50" = '  -3.40%  '
    "D51" = '6.990'
    "E51" = '  -5.91%  '
}

foreach ($addr in $updates.Keys) {
    Set-TextValue $ws $addr $updates[$addr]
}

$excel.CutCopyMode = 0